# Main_Data.xlsx edit — "Add files via upload"
#
# Net effect (derived from the OOXML diff):
#  1. workbook.xml: absPath tweaked to include "/Desktop/"; the "Overall"
#     sheet becomes the active/selected tab instead of "ALL".
#  2. Overall sheet ("Overall" = sheet2):
#       - The three federal-territory rows WPKL / WPLABUAN / WPPUTRAJAYA
#         are merged into a single "WPERSEKUTUAN" row (row 4), whose
#         2015/2016/2017 totals are now live SUM-style formulas.
#       - The OVERSEA row is moved up to row 7 (right after SELANGOR),
#         and KEDAH..PERLIS shift up to fill rows 8-16.
#       - Rows 17-18 (now empty) are removed; column A is widened to fit
#         "WPERSEKUTUAN"; selection moves to B4.
#     (The shared-strings table shrinks by 3 unused entries and grows by
#     the new "WPERSEKUTUAN" entry purely as a side effect of the above —
#     no script action needed for that.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

# --- Row 4: WPKL + WPLABUAN + WPPUTRAJAYA -> WPERSEKUTUAN (formulas) ---
$ws.Range("A4").Value = "WPERSEKUTUAN"
$ws.Range("B4").Formula = "=11423+202+27"
$ws.Range("C4").Formula = "=13465+232+38"
$ws.Range("D4").Formula = "=11060+240+105"

# --- Rows 5-6 (PERAK, SELANGOR) are unchanged ---

# --- Row 7: OVERSEA moves up here; matches the old OVERSEA row's bare
#     (unstyled) B/D formatting, with C keeping the #,##0 style ---
$ws.Range("A7").Value = "OVERSEA "
$ws.Range("B7:D7").ClearFormats() | Out-Null
$ws.Range("B7").Value = 458
$ws.Range("C7").Value = 495
$ws.Range("C7").NumberFormat = "#,##0"
$ws.Range("D7").Value = 347

# --- Rows 8-16: KEDAH..PERLIS shifted up one row (same #,##0 styling
#     the rows already carried one slot down) ---
$ws.Range("A8").Value = "KEDAH "
$ws.Range("B8").Value = 9740
$ws.Range("C8").Value = 10715
$ws.Range("D8").Value = 9495

$ws.Range("A9").Value = "TERENGGANU"
$ws.Range("B9").Value = 7403
$ws.Range("C9").Value = 8662
$ws.Range("D9").Value = 7380

$ws.Range("A10").Value = "PAHANG"
$ws.Range("B10").Value = 7270
$ws.Range("C10").Value = 8195
$ws.Range("D10").Value = 6933

$ws.Range("A11").Value = "SARAWAK "
$ws.Range("B11").Value = 7259
$ws.Range("C11").Value = 8036
$ws.Range("D11").Value = 6520

$ws.Range("A12").Value = "SABAH "
$ws.Range("B12").Value = 6531
$ws.Range("C12").Value = 7851
$ws.Range("D12").Value = 6872

$ws.Range("A13").Value = "PULAUPINANG "
$ws.Range("B13").Value = 6111
$ws.Range("C13").Value = 6740
$ws.Range("D13").Value = 5789

$ws.Range("A14").Value = "NSEMBILAN "
$ws.Range("B14").Value = 4285
$ws.Range("C14").Value = 4865
$ws.Range("D14").Value = 4287

$ws.Range("A15").Value = "MELAKA "
$ws.Range("B15").Value = 3988
$ws.Range("C15").Value = 4399
$ws.Range("D15").Value = 3717

# Row 16: PERLIS values move up from row 15, and pick up the s="12"/s="8"
# numeric styling (the old row16/17/18 bare cells are being retired below).
$ws.Range("A16").Value = "PERLIS "
$ws.Range("B16:D16").ClearFormats() | Out-Null
$ws.Range("B16").Value = 1343
$ws.Range("B16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 1631
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 1429
$ws.Range("D16").NumberFormat = "#,##0"

# --- Old rows 17 & 18 (WPPUTRAJAYA, OVERSEA) no longer needed ---
$ws.Rows.Item(17).Delete() | Out-Null
$ws.Rows.Item(17).Delete() | Out-Null

# --- Column A widens to fit the new longest label ("WPERSEKUTUAN") ---
$ws.Columns.Item(1).AutoFit() | Out-Null

# --- View state: Overall becomes the active sheet/tab, selection on B4 ---
$ws.Activate()
$ws.Range("B4").Select() | Out-Null

